# Insert a new data row at row 469 (pushing existing rows 469-534 down to
# 470-535), then populate the new row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("469:469").Insert()

$ws.Cells.Item(469, 1).Value  = 11
$ws.Cells.Item(469, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(469, 3).Value  = "Bíobío"
$ws.Cells.Item(469, 4).Value  = 45127
$ws.Cells.Item(469, 5).Value  = 8
$ws.Cells.Item(469, 6).Value  = 100114014
$ws.Cells.Item(469, 7).Value  = "Betarraga"
$ws.Cells.Item(469, 8).Value  = "Sin especificar"
$ws.Cells.Item(469, 9).Value  = "Primera"
$ws.Cells.Item(469, 10).Value = 700
$ws.Cells.Item(469, 11).Value = 650
$ws.Cells.Item(469, 12).Value = 700
$ws.Cells.Item(469, 13).Value = 679
$ws.Cells.Item(469, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(469, 15).Value = "Región Metropolitana"
$ws.Cells.Item(469, 16).Value = 136
$ws.Cells.Item(469, 17).Value = 5
$ws.Cells.Item(469, 18).Value = "Hortaliza"
